$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BuscaBarraFalha")
$ws.Range("E5").Style = "Comma"
$ws.Range("E5").HorizontalAlignment = -4108
